$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.8564554931836408
$ws.Range("C2").Value = 0.9393139841688655
$ws.Range("D2").Value = 0.8959731543624161
$ws.Range("E2").Value = 1137

$ws.Range("B3").Value = 0.9542586750788643
$ws.Range("C3").Value = 0.9438377535101404
$ws.Range("D3").Value = 0.9490196078431372
$ws.Range("E3").Value = 641

$ws.Range("B4").Value = 0.8596713021491783
$ws.Range("C4").Value = 0.8436724565756824
$ws.Range("D4").Value = 0.8515967438948028
$ws.Range("E4").Value = 806

$ws.Range("B5").Value = 0.8544776119402985
$ws.Range("C5").Value = 0.6432584269662921
$ws.Range("D5").Value = 0.7339743589743589
$ws.Range("E5").Value = 356

$ws.Range("B6").Value = 0.8782312925170068
$ws.Range("C6").Value = 0.8782312925170068
$ws.Range("D6").Value = 0.8782312925170068
$ws.Range("E6").Value = 0.8782312925170068

$ws.Range("B7").Value = 0.8812157705879955
$ws.Range("C7").Value = 0.8425206553052451
$ws.Range("D7").Value = 0.8576409662686788

$ws.Range("B8").Value = 0.8784213625368488
$ws.Range("C8").Value = 0.8782312925170068
$ws.Range("D8").Value = 0.8757567661604083
